$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Target values for rows 31-43 (row 38 unchanged) across columns A,B,D,E,F,G,H,I,J,Q,R
# Derived from permuting the existing row data per the source diff.

$data = @{
    31 = @{ A=112017130; B=90814;  D="LC"; E=4364; F="Dropptaggsvamp";        G="Hydnellum ferrugineum"; H="(Fr.:Fr.) P. Karst.";      I=$null; J=$null;          Q=682695; R=6575454 }
    32 = @{ A=112017326; B=90808;  D="NT"; E=4362; F="Blå taggsvamp";         G="Hydnellum caeruleum";   H="(Hornem.) P.Karst.";        I=$null; J=$null;          Q=682714; R=6575496 }
    33 = @{ A=112017447; B=90814;  D="LC"; E=4364; F="Dropptaggsvamp";        G="Hydnellum ferrugineum"; H="(Fr.:Fr.) P. Karst.";      I=$null; J=$null;          Q=682844; R=6575514 }
    34 = @{ A=112017488; B=90826;  D="LC"; E=4366; F="Skarp dropptaggsvamp";  G="Hydnellum peckii";      H="Banker";                    I=$null; J=$null;          Q=682956; R=6575474 }
    35 = @{ A=112017512; B=88180;  D="VU"; E=6276; F="Goliatmusseron";        G="Tricholoma matsutake";  H="(S.Ito & S.Imai) Singer";   I=4;     J="fruktkroppar"; Q=683037; R=6575484 }
    36 = @{ A=112017413; B=90857;  D="NT"; E=5448; F="Svartvit taggsvamp";    G="Phellodon connatus";    H="(Schultz) nom.prov";        I=$null; J=$null;          Q=682734; R=6575482 }
    37 = @{ A=112017252; B=90814;  D="LC"; E=4364; F="Dropptaggsvamp";        G="Hydnellum ferrugineum"; H="(Fr.:Fr.) P. Karst.";      I=$null; J=$null;          Q=682711; R=6575494 }
    39 = @{ A=112017534; B=88140;  D="VU"; E=1593; F="Lakritsmusseron";       G="Tricholoma apium";      H="Jul.Schäff.";               I=4;     J="fruktkroppar"; Q=683073; R=6575478 }
    40 = @{ A=112017392; B=90858;  D="NT"; E=5449; F="Svart taggsvamp";       G="Phellodon niger";       H="(Fr.:Fr.) P.Karst.";        I=$null; J=$null;          Q=682712; R=6575458 }
    41 = @{ A=112017465; B=88180;  D="VU"; E=6276; F="Goliatmusseron";        G="Tricholoma matsutake";  H="(S.Ito & S.Imai) Singer";   I=3;     J="fruktkroppar"; Q=682896; R=6575514 }
    42 = @{ A=112017224; B=90826;  D="LC"; E=4366; F="Skarp dropptaggsvamp";  G="Hydnellum peckii";      H="Banker";                    I=$null; J=$null;          Q=682703; R=6575491 }
    43 = @{ A=112017430; B=90857;  D="NT"; E=5448; F="Svartvit taggsvamp";    G="Phellodon connatus";    H="(Schultz) nom.prov";        I=$null; J=$null;          Q=682793; R=6575520 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
}
